$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume(1h) (E) values per the Feb 9 2023 GitHub Actions refresh commit.
# Leading apostrophe forces these numeric/percent-looking strings to stay stored as literal text,
# matching the source data (originally inline strings, e.g. "317.77", "-3.28%").
$ws.Range("D2").Value = "'317.95"
$ws.Range("E2").Value = "'-3.22%"

$ws.Range("D3").Value = "'42.04"
$ws.Range("E3").Value = "'-5.13%"

$ws.Range("D4").Value = "'5.192"
$ws.Range("E4").Value = "'2.51%"

$ws.Range("D5").Value = "'0.08110"
$ws.Range("E5").Value = "'-3.15%"

$ws.Range("D6").Value = "'4.374"
$ws.Range("E6").Value = "'-1.59%"

$ws.Range("D7").Value = "'1.753"
$ws.Range("E7").Value = "'-10.33%"

$ws.Range("D8").Value = "'0.9308"
$ws.Range("E8").Value = "'-4.70%"

$ws.Range("D9").Value = "'0.1121"
$ws.Range("E9").Value = "'-1.55%"

$ws.Range("D10").Value = "'0.1856"
$ws.Range("E10").Value = "'-2.53%"

$ws.Range("D11").Value = "'0.09305"

$ws.Range("D12").Value = "'0.04565"
$ws.Range("E12").Value = "'-1.51%"

$ws.Range("D13").Value = "'7.395"
$ws.Range("E13").Value = "'-16.92%"

$ws.Range("E14").Value = "'-0.55%"

$ws.Range("D15").Value = "'0.001296"
$ws.Range("E15").Value = "'-0.08%"

$ws.Range("D16").Value = "'0.005988"
$ws.Range("E16").Value = "'-0.40%"

$ws.Range("D17").Value = "'3.361"
$ws.Range("E17").Value = "'-1.21%"

$ws.Range("D18").Value = "'2.588"
$ws.Range("E18").Value = "'3.48%"

$ws.Range("D19").Value = "'0.3356"
$ws.Range("E19").Value = "'1.03%"

$ws.Range("D20").Value = "'0.1382"
$ws.Range("E20").Value = "'2.10%"

$ws.Range("D21").Value = "'0.2549"
$ws.Range("E21").Value = "'-0.14%"

$ws.Range("D22").Value = "'0.04180"
$ws.Range("E22").Value = "'0.42%"

$ws.Range("D23").Value = "'0.001243"
$ws.Range("E23").Value = "'-4.26%"

$ws.Range("D24").Value = "'0.004289"
$ws.Range("E24").Value = "'-2.59%"

$ws.Range("E25").Value = "'-6.21%"

$ws.Range("D26").Value = "'0.0002983"
$ws.Range("E26").Value = "'-0.07%"

$ws.Range("D38").Value = "'0.02587"
$ws.Range("E38").Value = "'-5.63%"

$ws.Range("D39").Value = "'0.05483"
$ws.Range("E39").Value = "'-2.73%"

$ws.Range("D40").Value = "'0.008038"
$ws.Range("E40").Value = "'1.86%"

$ws.Range("D41").Value = "'0.1393"
$ws.Range("E41").Value = "'-1.53%"

$ws.Range("D42").Value = "'0.007366"
$ws.Range("E42").Value = "'0.00%"

$ws.Range("D43").Value = "'0.002090"
$ws.Range("E43").Value = "'-1.11%"

$ws.Range("D44").Value = "'0.008255"
$ws.Range("E44").Value = "'4.22%"

$ws.Range("D45").Value = "'0.3450"
$ws.Range("E45").Value = "'-1.68%"

$ws.Range("D46").Value = "'0.00006734"
$ws.Range("E46").Value = "'-2.56%"

$ws.Range("E47").Value = "'-0.07%"

$ws.Range("D48").Value = "'0.003375"
$ws.Range("E48").Value = "'-3.97%"

$ws.Range("D49").Value = "'0.004106"
$ws.Range("E49").Value = "'16.07%"

$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'-0.07%"

$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'-0.07%"

